$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (F1:L1)
$ws.Range("F1").Value = "d37"
$ws.Range("G1").Value = "i74"
$ws.Range("H1").Value = "bl97"
$ws.Range("I1").Value = "v94"
$ws.Range("J1").Value = "b12+"
$ws.Range("K1").Value = "b65"
$ws.Range("L1").Value = "z20"

# Row 2
$ws.Range("F2").Value = 357654105415049
$ws.Range("G2").Value = 359999101446706
$ws.Range("H2").Value = 357655105029508
$ws.Range("I2").Value = 359311090420736
$ws.Range("J2").Value = 353428112075748
$ws.Range("K2").Value = 357653106490787
$ws.Range("L2").Value = 358126101353648

# Row 3
$ws.Range("F3").Value = 357654105295508
$ws.Range("G3").Value = 359999101441921
$ws.Range("H3").Value = 357655105031264
$ws.Range("J3").Value = 353428112075763
$ws.Range("K3").Value = 357653106252666
$ws.Range("L3").Value = 358126101354885

# Row 4
$ws.Range("F4").Value = 357654105273729
$ws.Range("G4").Value = 359999101427581
$ws.Range("H4").Value = 357655105083760
$ws.Range("J4").Value = 353428112000308
$ws.Range("K4").Value = 357653106248763

# Row 5
$ws.Range("F5").Value = 357654105274685
$ws.Range("H5").Value = 357655105087084
$ws.Range("J5").Value = 353428112074220
$ws.Range("K5").Value = 357653104936500

# Row 6
$ws.Range("F6").Value = 357654105273620
$ws.Range("J6").Value = 353428112080102
$ws.Range("K6").Value = 357653106651727

# Row 7
$ws.Range("F7").Value = 357654105274040
$ws.Range("J7").Value = 353428110202302
$ws.Range("K7").Value = 357653106486421

# Row 8
$ws.Range("F8").Value = 357654105280187
$ws.Range("J8").Value = 353428112000282
$ws.Range("K8").Value = 357653106251304

# Row 9
$ws.Range("F9").Value = 357654105280427
$ws.Range("J9").Value = 353428110205487

# Row 10
$ws.Range("J10").Value = 353428110227804

# Row 11
$ws.Range("J11").Value = 353428112063983

# Row 12
$ws.Range("J12").Value = 353428112080128

# Row 13
$ws.Range("J13").Value = 353428112000266

# Row 14
$ws.Range("J14").Value = 353428112075789

# Row 15
$ws.Range("J15").Value = 353428112087560

# The new columns (F:L) hold 15-digit IMEI numbers just like the existing
# "best fit" columns C:E, so widen them to the same best-fit width.
$ws.Columns("F:L").ColumnWidth = 15.25

# Update selection to match the final state of the diff
$ws.Range("H21").Select()
